$wb = $excel.ActiveWorkbook

# --- Rename the "Include from LOINC" sheet to "Include #0" ---
$ws2 = $wb.Worksheets.Item("Include from LOINC")
$ws2.Name = "Include #0"

# --- Update the "Metadata" sheet ---
$ws1 = $wb.Worksheets.Item("Metadata")

# Insert a new "Jurisdiction" property row above the existing "Description" row
# (row 11), pushing Description/Purpose/Copyright/Immutable down by one row.
# Use Copy(destination) (bottom-up) so the existing cell style (s="2") carries
# over to the newly used rows instead of picking up a freshly-minted default
# style.
$ws1.Range("A14:B14").Copy($ws1.Range("A15:B15"))
$ws1.Range("A13:B13").Copy($ws1.Range("A14:B14"))
$ws1.Range("A12:B12").Copy($ws1.Range("A13:B13"))
$ws1.Range("A11:B11").Copy($ws1.Range("A12:B12"))

# Copy() only carries the style forward when the source cell already holds a
# value; a blank source leaves the destination's old value in place, so the
# now-blank "Purpose" value cell (B13) needs to be cleared explicitly.
$ws1.Range("B13").Value = ""

# Populate the new Jurisdiction row.
$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = ""

# Refresh the publication Date value.
$ws1.Range("B8").Value = "2024-09-17T19:55:11+00:00"
